# Add new user/sucursal records to the "Users" sheet.
# (Commit: "Se agregaron casos en los siguientes modulos
#  ASUC58 - ASUC59 , MEP04 - MEP05 - MEP06 , PAS12")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Helper to write a "Sucursal" (column C) cell as text, preserving any
# leading zeros and matching the right-aligned text styling used by the
# rest of the column.
function Set-Sucursal($rowNum, $value) {
    $cCell = $ws.Cells.Item($rowNum, 3)
    $cCell.NumberFormat = "@"
    $cCell.HorizontalAlignment = -4152
    $cCell.Value = $value
}

# Populate the new rows. The cells are written in the same order the
# original data entry happened in, so that newly introduced shared
# strings land at the same indices as the source workbook.
$ws.Cells.Item(60, 1).Value = "F00688"

$ws.Cells.Item(57, 1).Value = "F00020"
$ws.Cells.Item(58, 1).Value = "F00847"
$ws.Cells.Item(59, 1).Value = "F03808"
Set-Sucursal 59 "322"

$ws.Cells.Item(61, 1).Value = "F00460"
Set-Sucursal 61 "060"

# Remaining "Sucursal" values reuse strings already present in the
# workbook, so write order no longer matters for these.
Set-Sucursal 57 "020"
Set-Sucursal 58 "074"
Set-Sucursal 60 "089"

# Move the active selection to just past the newly added data, like the
# original author left the sheet positioned at C62 after editing C61.
$ws.Range("C62").Select() | Out-Null
